$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "curva" (A) and "fecha" (B) columns, shifting "plazo" (C) and
# "tasa" (D) left into A and B.
$ws.Range("A:B").Delete()

# Match the saved selection state (cell C8, which maps onto the "tasa"
# column after the shift).
[void]$ws.Range("C8").Select()
